# Update data-example to 32-byte pubkey
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A header: "Mã giáo vụ" -> "Mã giảng viên"
$ws.Range("A1").Value = "Mã giảng viên"

# Column E (Khóa công khai) values -> full 32-byte (64 hex char) pubkeys prefixed with 02/03
$ws.Range("E2").Value = "02a2e3f2b6b9ec1155979ee691072bd17fc9facd0d7751417fdf4d7af760ce962a"
$ws.Range("E3").Value = "02992cf23456bf4167fc2d69f70ee36361af8ad2a2512f9a660be320086211f3e0"
$ws.Range("E4").Value = "028b921e7d9e7e0ece5660dc48e6c7b84d8c547a3f3d42990e67897858968120d5"
$ws.Range("E5").Value = "02a09d011afd12783ab87f44ad1932dcca1f59e13ba904d2c7b465df31df681a32"
$ws.Range("E6").Value = "0351befccfa9fa6ba05d16a2eb94a593f773f7e9db0138ba831424905174d730fe"

# Match final selection state observed in the saved file
$ws.Range("E6").Select()
